$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "Total Pages Read"
$ws.Range("B11").Value = 1402
